$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-ThinBox($rng) {
    $b = $rng.Borders
    $b.Item(7).LineStyle = 1
    $b.Item(10).LineStyle = 1
    $b.Item(7).Weight = 2
    $b.Item(10).Weight = 2
    $b.Item(7).ColorIndex = 1
    $b.Item(10).ColorIndex = 1
}

# Row 217 - หมูเด้ง (Moo Deng) live cam, Chonburi, Thailand
$ws.Range("F217").Value = "7EEy1OEmGjc"
$ws.Range("B217").Value = "13.214975102821438, 101.0569795387831"
$ws.Range("C217").Value = "หมูเด้ง - Moo Deng Live 🔴"
$ws.Range("D217").Value = "Chonburi"
$ws.Range("A217").Value = "LIVE, ZOO"
$ws.Range("E217").Value = "Thailand"
Set-ThinBox $ws.Range("A217")
Set-ThinBox $ws.Range("E217")

# Row 218 - Seibu Ikebukuro Line camera 2 (electric storage line), Tokyo, Japan
$ws.Range("F218").Value = "kzQdszcQ2HM"
$ws.Range("B218").Value = "35.748233007275104, 139.56447498450478"
$ws.Range("C218").Value = "【西武線ライブカメラ２】　西武池袋線保谷駅付近の電留線ライブ映像/Live footage of the electric storage line near Hoya Station on the Seibu Ikebukuro Line"
$ws.Range("D218").Value = "Tokyo"
$ws.Range("A218").Value = "LIVE, TRAIN, RAIL"
$ws.Range("E218").Value = "Japan"
Set-ThinBox $ws.Range("A218")
Set-ThinBox $ws.Range("E218")

# Row 219 - Seibu Ikebukuro Line camera (railroad crossing), Tokyo, Japan (no Category)
$ws.Range("B219").Value = "35.74805277396905, 139.56531155077093"
$ws.Range("F219").Value = "TfzJPS1rJmk"
$ws.Range("C219").Value = "【西武線ライブカメラ】　西武池袋線保谷駅付近の踏切からライブ配信/Live broadcast from a railroad crossing near Hoya Station on the Seibu Ikebukuro Line"
$ws.Range("D219").Value = "Tokyo"
$ws.Range("E219").Value = "Japan"
Set-ThinBox $ws.Range("E219")

$ws.Range("C219").Select()
